# Update Receptor/Edge specificity columns (M:T) on the single worksheet
# with refreshed TPM-derived values, per "update scripts wuth new tpm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.05871566666666667
$ws.Range("N2").Value = 0.176147
$ws.Range("O2").Value = 0.008355169877275808
$ws.Range("P2").Value = 0.008355169877275808
$ws.Range("Q2").Value = 0.08616567141488887
$ws.Range("R2").Value = 0.7754910427339999
$ws.Range("S2").Value = 0.001579757496688954
$ws.Range("T2").Value = 0.001579757496688954

# Row 3
$ws.Range("O3").Value = 0.1868088427899751
$ws.Range("P3").Value = 0.1868088427899751
$ws.Range("S3").Value = 0.03532096584270439
$ws.Range("T3").Value = 0.03532096584270439

# Row 4
$ws.Range("O4").Value = 0.8048359873327491
$ws.Range("P4").Value = 0.8048359873327491
$ws.Range("S4").Value = 0.1521747257410067
$ws.Range("T4").Value = 0.1521747257410067

# Row 5
$ws.Range("M5").Value = 0.05871566666666667
$ws.Range("N5").Value = 0.176147
$ws.Range("O5").Value = 0.008355169877275808
$ws.Range("P5").Value = 0.008355169877275808
$ws.Range("Q5").Value = 0.3695554274055555
$ws.Range("R5").Value = 3.32599884665
$ws.Range("S5").Value = 0.006775412380586854
$ws.Range("T5").Value = 0.006775412380586853

# Row 6
$ws.Range("O6").Value = 0.1868088427899751
$ws.Range("P6").Value = 0.1868088427899751
$ws.Range("S6").Value = 0.1514878769472707
$ws.Range("T6").Value = 0.1514878769472707

# Row 7
$ws.Range("O7").Value = 0.8048359873327491
$ws.Range("P7").Value = 0.8048359873327491
$ws.Range("S7").Value = 0.6526612615917424
$ws.Range("T7").Value = 0.6526612615917423
